$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.102735638618469
$ws.Range("B1").Value = 1.945701360702515
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.008204460144043
$ws.Range("E1").Value = 1.120494723320007
